$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Text Translated")

# Row 1 header language code updates
$ws.Range("D1").Value = "fr"
$ws.Range("E1").Value = "de"
$ws.Range("F1").Value = "el"

# Row 2 content updates
$ws.Range("A2").Value = "Attn"
$ws.Range("B2").Value = "FFF2CC"
$ws.Range("C2").Value = "注意"
$ws.Range("D2").Value = "À l'attention de"
$ws.Range("E2").Value = "Beachtung"
$ws.Range("F2").Value = "Προσοχή"
$ws.Range("G2").Value = "Attenzione"
$ws.Range("H2").Value = "주목"
$ws.Range("I2").Value = "Uwaga"
$ws.Range("J2").Value = "Atenção"
$ws.Range("K2").Value = "Atenção"
$ws.Range("L2").Value = "Atención"
$ws.Range("M2").Value = "Atención"
$ws.Range("N2").Value = "ความสนใจ"
$ws.Range("O2").Value = "Liên Hệ"
$ws.Range("P2").Value = "CustomLabel`$CEC_Attention"
$ws.Range("Q2").Value = "CustomLabel"
$ws.Range("R2").Value = "Attn"
$ws.Range("S2").Value = "IVP Contact Creation Translations.xlsx"
$ws.Range("T2").Value = "Sheet1"
